# Update to the Jul-Dec 2022 reporting period (3 de marzo 2023 update)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- Row 8: reporting period dates ---
$ws.Range("B8").Value = "7/1/2022"
$ws.Range("C8").Value = "12/31/2022"
$ws.Range("U8").Value = "1/10/2023"
$ws.Range("V8").Value = "1/10/2023"

# --- T8: área responsable text changed ---
$ws.Range("T8").Value = "Subdirección de Recursos Financieros (UPP)"

# --- Row heights ---
$ws.Rows.Item(3).RowHeight = 37.5
$ws.Rows.Item(8).RowHeight = 123

# --- Column W width ---
$ws.Columns.Item(23).ColumnWidth = 43.140625

# --- Data validation ranges shrink from row 201 to row 124 ---
$ws.Range("D8:D201").Validation.Delete()
$ws.Range("D8:D124").Validation.Add(3, 1, 1, "=Hidden_13")
$ws.Range("R8:R201").Validation.Delete()
$ws.Range("R8:R124").Validation.Add(3, 1, 1, "=Hidden_217")

# --- Selection moved ---
$ws.Range("B11").Select()

# --- Cell format tweaks: T8 loses its highlight style, becomes like the rest of the row ---
$ws.Range("T8").Interior.Pattern = -4142
$ws.Range("T8").HorizontalAlignment = -4131
$ws.Range("T8").Font.Name = "Calibri"

# --- W8 keeps justify+wrap but loses its gray fill ---
$ws.Range("W8").Interior.Pattern = -4142
$ws.Range("W8").HorizontalAlignment = -4130
$ws.Range("W8").WrapText = $true

# --- G3 header cell border simplified to left-only ---
$ws.Range("G3").Borders.LineStyle = 0
$ws.Range("G3").Borders.Item(7).LineStyle = 1
$ws.Range("G3").Borders.Item(7).Weight = 2
